$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts existing column B -> C, which
# keeps its original 50.21875 width automatically)
$ws.Columns.Item(2).Insert()

# Set the new column B to its target width (15.1 "set" width renders as
# the quantized 16-character stored width seen in the target file)
$ws.Columns.Item(2).ColumnWidth = 15.1

# Apply bold style (same as column A header cells) to the new B cells in rows 1-3,5-7 (kept empty)
$ws.Range("B1:B3").Font.Bold = $true
$ws.Range("B5:B7").Font.Bold = $true

# Fill in the new code values in column B for rows 8-12
# (entry order mirrors shared-string insertion order of the original edit)
$ws.Range("B9").Value = "A05"
$ws.Range("B8").Value = "A08"
$ws.Range("B10").Value = "C02"
$ws.Range("B11").Value = "C04, C07"
$ws.Range("B12").Value = "R1"

# Add the new row 12 description in column C
$ws.Range("C12").Value = "Requirements ändern um sie an der Anwendung besser anzupassen "

# Update selection to match the final state
$ws.Range("A11").Select()
